$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new H1 header cell
# so the new "Save" column header matches the other header cells' style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
